$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at position 174 (a new weekly price observation),
# pushing the existing rows 174..256 down to 175..257.
$ws.Rows.Item(174).Insert()

# Populate the newly inserted row 174 with the new observation. Every
# column except Fecha/Precio min/max/promedio/Precio $/Kg matches the
# constant metadata shared by every row in this block.
$ws.Cells.Item(174, 1).Value = 7
$ws.Cells.Item(174, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(174, 3).Value = "Ñuble"
$ws.Cells.Item(174, 4).Value = 44609
$ws.Cells.Item(174, 5).Value = 16
$ws.Cells.Item(174, 6).Value = 100114013
$ws.Cells.Item(174, 7).Value = "Zanahoria"
$ws.Cells.Item(174, 8).Value = "Sin especificar"
$ws.Cells.Item(174, 9).Value = "Primera"
$ws.Cells.Item(174, 10).Value = 120
$ws.Cells.Item(174, 11).Value = 7000
$ws.Cells.Item(174, 12).Value = 7500
$ws.Cells.Item(174, 13).Value = 7250
$ws.Cells.Item(174, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(174, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(174, 16).Value = 362
$ws.Cells.Item(174, 17).Value = 20
$ws.Cells.Item(174, 18).Value = "Hortaliza"
